# Update 20200421 - incluye Recuperados Cba
#
# 1. Add 9 retroactive "recuperados" (E column) values for Córdoba on
#    existing rows in "reporte_vespertino".
# 2. Append a new day (2020-04-21) of rows to "reporte_vespertino" for
#    every province/jurisdiction.
# 3. Append the same new day's two rows to "reporte_matutino".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Historical corrections: Cordoba "recuperados" (column E) on
#    reporte_vespertino for several already-existing rows.
# ---------------------------------------------------------------------
$wsVesp = $wb.Worksheets.Item("reporte_vespertino")

$corrections = @(
    @(734, 1),
    @(909, 13),
    @(934, 17),
    @(959, 2),
    @(984, 4),
    @(1059, 1),
    @(1158, 5),
    @(1183, 2),
    @(1208, 4)
)

foreach ($fix in $corrections) {
    $r = $fix[0]
    $val = $fix[1]
    $wsVesp.Cells.Item($r, 5).Value = $val
}

# ---------------------------------------------------------------------
# 2) New rows for 2020-04-21 on reporte_vespertino (rows 1227-1251)
# ---------------------------------------------------------------------
$newRowsVesp = @(
    @(1227, "Argentina_Nacion", 112, 6, 0),
    @(1228, "Buenos Aires", 60, 1, 0),
    @(1229, "CABA", 33, 3, 0),
    @(1230, "Catamarca", 0, 0, 0),
    @(1231, "Chaco", 5, 0, 6),
    @(1232, "Chubut", 0, 0, 0),
    @(1233, "Córdoba", 0, 0, 9),
    @(1234, "Corrientes", 0, 0, 0),
    @(1235, "Entre Ríos", 0, 0, 0),
    @(1236, "Formosa", 0, 0, 0),
    @(1237, "Jujuy", 0, 0, 0),
    @(1238, "La Pampa", 0, 0, 0),
    @(1239, "La Rioja", 1, 2, 0),
    @(1240, "Mendoza", 0, 0, 0),
    @(1241, "Misiones", 0, 0, 0),
    @(1242, "Neuquén", 0, 0, 0),
    @(1243, "Rio Negro", 11, 0, 0),
    @(1244, "Salta", 0, 0, 0),
    @(1245, "San Juan", 0, 0, 0),
    @(1246, "San Luis", 0, 0, 0),
    @(1247, "Santa Cruz", 0, 0, 0),
    @(1248, "Santa Fe", 1, 0, 12),
    @(1249, "Santiago del Estero", 0, 0, 0),
    @(1250, "Tierra del Fuego", 1, 0, 0),
    @(1251, "Tucumán", 0, 0, 0)
)

$lastRowVesp = 1227 + $newRowsVesp.Count - 1
$wsVesp.Range("A1227:A$lastRowVesp").NumberFormat = "@"

foreach ($row in $newRowsVesp) {
    $r = $row[0]
    $wsVesp.Cells.Item($r, 1).Value = "2020-04-21"
    $wsVesp.Cells.Item($r, 2).Value = $row[1]
    $wsVesp.Cells.Item($r, 3).Value = $row[2]
    $wsVesp.Cells.Item($r, 4).Value = $row[3]
    $wsVesp.Cells.Item($r, 5).Value = $row[4]
}

# Scroll reporte_vespertino back to the top (matches topLeftCell A1).
$wsVesp.Activate()
$wsVesp.Range("A1").Select()

# ---------------------------------------------------------------------
# 3) New rows for 2020-04-21 on reporte_matutino (rows 38-39)
# ---------------------------------------------------------------------
$wsMat = $wb.Worksheets.Item("reporte_matutino")

$newRowsMat = @(
    @(38, "Argentina_Nacion", 0, 3, 103),
    @(39, "Buenos Aires", 0, 3, 0)
)

$lastRowMat = 38 + $newRowsMat.Count - 1
$wsMat.Range("A38:A$lastRowMat").NumberFormat = "@"

foreach ($row in $newRowsMat) {
    $r = $row[0]
    $wsMat.Cells.Item($r, 1).Value = "2020-04-21"
    $wsMat.Cells.Item($r, 2).Value = $row[1]
    $wsMat.Cells.Item($r, 3).Value = $row[2]
    $wsMat.Cells.Item($r, 4).Value = $row[3]
    $wsMat.Cells.Item($r, 5).Value = $row[4]
}

$wsMat.Activate()
$wsMat.Range("E$lastRowMat").Select()

# ---------------------------------------------------------------------
# Restore reporte_vespertino as the active sheet/tab, as in the source.
# ---------------------------------------------------------------------
$wsVesp.Activate()
$wsVesp.Range("A1").Select()
